# "Fruta / hortaliza, semanal"
# The weekly refresh reshuffled the Fecha (D), Volumen (J), Precio minimo (K),
# Precio maximo (L), Precio promedio ponderado (M) and Precio $/Kg (P) values
# across the existing data rows (2-20). All other columns are unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row => Fecha, Volumen, Precio minimo, Precio maximo, Precio promedio ponderado, Precio $/Kg
$rows = @(
    @{ Row = 2;  D = 44517; J = 500;  K = 800;  L = 900;  M = 850;  P = 850 }
    @{ Row = 3;  D = 44525; J = 360;  K = 800;  L = 900;  M = 850;  P = 850 }
    @{ Row = 4;  D = 44503; J = 400;  K = 900;  L = 1000; M = 950;  P = 950 }
    @{ Row = 5;  D = 44537; J = 400;  K = 800;  L = 900;  M = 850;  P = 850 }
    @{ Row = 6;  D = 44518; J = 400;  K = 800;  L = 900;  M = 850;  P = 850 }
    @{ Row = 7;  D = 44508; J = 400;  K = 900;  L = 1000; M = 950;  P = 950 }
    @{ Row = 8;  D = 44505; J = 440;  K = 900;  L = 1000; M = 950;  P = 950 }
    @{ Row = 9;  D = 44553; J = 8000; K = 800;  L = 900;  M = 850;  P = 850 }
    @{ Row = 10; D = 44476; J = 300;  K = 1100; L = 1200; M = 1150; P = 1150 }
    @{ Row = 11; D = 44524; J = 400;  K = 800;  L = 900;  M = 850;  P = 850 }
    @{ Row = 12; D = 44530; J = 300;  K = 800;  L = 900;  M = 850;  P = 850 }
    @{ Row = 13; D = 44510; J = 600;  K = 900;  L = 1000; M = 950;  P = 950 }
    @{ Row = 14; D = 44516; J = 400;  K = 900;  L = 1000; M = 950;  P = 950 }
    @{ Row = 15; D = 44512; J = 600;  K = 900;  L = 1000; M = 950;  P = 950 }
    @{ Row = 16; D = 44511; J = 500;  K = 900;  L = 1000; M = 950;  P = 950 }
    @{ Row = 17; D = 44545; J = 4000; K = 800;  L = 900;  M = 850;  P = 850 }
    @{ Row = 18; D = 44504; J = 500;  K = 900;  L = 1000; M = 950;  P = 950 }
    @{ Row = 19; D = 44532; J = 240;  K = 800;  L = 900;  M = 850;  P = 850 }
    @{ Row = 20; D = 44523; J = 400;  K = 800;  L = 900;  M = 850;  P = 850 }
)

foreach ($item in $rows) {
    $r = $item.Row
    $ws.Cells.Item($r, 4).Value  = $item.D   # D - Fecha
    $ws.Cells.Item($r, 10).Value = $item.J   # J - Volumen
    $ws.Cells.Item($r, 11).Value = $item.K   # K - Precio minimo
    $ws.Cells.Item($r, 12).Value = $item.L   # L - Precio maximo
    $ws.Cells.Item($r, 13).Value = $item.M   # M - Precio promedio ponderado
    $ws.Cells.Item($r, 16).Value = $item.P   # P - Precio $/Kg
}
